$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (% Recommend) keeps its original plain-text representation
# instead of being auto-converted to a percentage number by Excel type inference.
$pctRange = $ws.Range("D2:D129")
$pctRange.NumberFormat = "@"

$data = @(
    @(2, 'Le Wagon', 1, 1, '100.00%'),
    @(3, 'Base10 Academy', 1, 1, '100.00%'),
    @(4, 'Fullstack Academy', 6, 6, '100.00%'),
    @(5, 'Ruby On The Beach', 1, 1, '100.00%'),
    @(6, 'Science to Data Science', 1, 1, '100.00%'),
    @(7, 'The Data Incubator', 1, 1, '100.00%'),
    @(8, 'devCodeCamp', 4, 4, '100.00%'),
    @(9, 'Coder''s Lab', 1, 1, '100.00%'),
    @(10, 'DESIGNATION', 1, 1, '100.00%'),
    @(11, 'Code Institute', 1, 1, '100.00%'),
    @(12, 'AngelHack Education', 1, 1, '100.00%'),
    @(13, 'Betamore', 1, 1, '100.00%'),
    @(14, 'Academia de CÃ³digo', 4, 4, '100.00%'),
    @(15, 'CodeCraft School', 1, 1, '100.00%'),
    @(16, 'Stackademy', 1, 1, '100.00%'),
    @(17, 'Founders & Coders', 1, 1, '100.00%'),
    @(18, 'Omaha Code School', 2, 2, '100.00%'),
    @(19, 'Telegraph Academy', 2, 2, '100.00%'),
    @(20, 'Eleven Fifty Academy', 2, 2, '100.00%'),
    @(21, 'We Got Coders', 2, 2, '100.00%'),
    @(22, 'Grand Circus', 3, 3, '100.00%'),
    @(23, 'Dev League', 3, 3, '100.00%'),
    @(24, 'Code Platoon', 2, 2, '100.00%'),
    @(25, 'EHD Academy', 2, 2, '100.00%'),
    @(26, 'BoiseCodeWorks', 2, 2, '100.00%'),
    @(27, 'Coder Factory', 2, 2, '100.00%'),
    @(28, 'DigitalCrafts', 2, 2, '100.00%'),
    @(29, 'CODEcamp Charleston', 2, 2, '100.00%'),
    @(30, 'Big Nerd Ranch', 2, 2, '100.00%'),
    @(31, 'Codecademy Labs', 2, 2, '100.00%'),
    @(32, 'Grace Hopper Academy', 3, 3, '100.00%'),
    @(33, 'AcadGild', 2, 2, '100.00%'),
    @(34, 'Coder Foundry', 1, 1, '100.00%'),
    @(35, 'CodeMasters Academy', 1, 1, '100.00%'),
    @(36, 'Claim Academy', 1, 1, '100.00%'),
    @(37, 'Byte Academy', 1, 1, '100.00%'),
    @(38, 'Operation Spark', 1, 1, '100.00%'),
    @(39, 'Origin Code Academy', 1, 1, '100.00%'),
    @(40, 'Austin Coding Academy', 1, 1, '100.00%'),
    @(41, 'Code Union', 1, 1, '100.00%'),
    @(42, 'Code 42', 1, 1, '100.00%'),
    @(43, 'Ladies Learning Code', 1, 1, '100.00%'),
    @(44, 'Dev Academy', 13, 13, '100.00%'),
    @(45, 'Makers Academy', 12, 12, '100.00%'),
    @(46, 'SeedPaths', 1, 1, '100.00%'),
    @(47, 'LEARN Academy', 3, 3, '100.00%'),
    @(48, 'CodeaCamp', 8, 8, '100.00%'),
    @(49, 'Viking Code School', 10, 10, '100.00%'),
    @(50, 'Launch School (formerly Tealeaf Academy)', 1, 1, '100.00%'),
    @(51, 'Microsoft Research Data Science Summer School', 1, 1, '100.00%'),
    @(52, 'codeU', 1, 1, '100.00%'),
    @(53, 'Turing', 27, 26, '96.30%'),
    @(54, 'Free Code Camp is not a bootcamp - please scroll up and change answer to "no"', 15, 14, '93.33%'),
    @(55, 'Hack Reactor', 29, 27, '93.10%'),
    @(56, 'Flatiron School', 54, 50, '92.59%'),
    @(57, 'App Academy', 22, 20, '90.91%'),
    @(58, 'MakerSquare', 20, 18, '90.00%'),
    @(59, 'Ada', 9, 8, '88.89%'),
    @(60, 'Hackbright Academy', 22, 19, '86.36%'),
    @(61, 'Dev Bootcamp', 48, 41, '85.42%'),
    @(62, 'Prime Digital Academy', 30, 25, '83.33%'),
    @(63, 'Thinkful', 6, 5, '83.33%'),
    @(64, 'The Firehose Project', 16, 13, '81.25%'),
    @(65, 'DevMountain', 10, 8, '80.00%'),
    @(66, 'Nashville Software School', 5, 4, '80.00%'),
    @(67, 'CareerFoundry', 5, 4, '80.00%'),
    @(68, 'Software Guild', 10, 8, '80.00%'),
    @(69, 'Launch Academy', 10, 8, '80.00%'),
    @(70, 'Epicodus', 14, 11, '78.57%'),
    @(71, 'General Assembly', 90, 70, '77.78%'),
    @(72, 'The Iron Yard', 40, 31, '77.50%'),
    @(73, 'LAMP Camp', 4, 3, '75.00%'),
    @(74, 'Skillcrush', 4, 3, '75.00%'),
    @(75, 'Lighthouse Labs', 4, 3, '75.00%'),
    @(76, 'Rutgers Coding Bootcamp', 4, 3, '75.00%'),
    @(77, 'Sabio.la', 4, 3, '75.00%'),
    @(78, 'Anyone Can Learn To Code', 7, 5, '71.43%'),
    @(79, 'IronHack', 7, 5, '71.43%'),
    @(80, 'CodeNinja', 7, 5, '71.43%'),
    @(81, 'Code Fellows', 21, 15, '71.43%'),
    @(82, 'New York Code + Design Academy', 20, 14, '70.00%'),
    @(83, 'Learn.Modern-Developer', 3, 2, '66.67%'),
    @(84, 'Tech Academy Portland', 3, 2, '66.67%'),
    @(85, 'Tech Talent South', 12, 8, '66.67%'),
    @(86, 'Orange County Code School', 3, 2, '66.67%'),
    @(87, 'BrainStation', 3, 2, '66.67%'),
    @(88, 'Coding Dojo', 18, 11, '61.11%'),
    @(89, 'Codeup', 10, 6, '60.00%'),
    @(90, 'Bitmaker Labs', 5, 3, '60.00%'),
    @(91, '10x.org.il', 5, 3, '60.00%'),
    @(92, 'CodeCore Bootcamp', 5, 3, '60.00%'),
    @(93, 'Startup Institute', 7, 4, '57.14%'),
    @(94, 'Galvanize', 20, 11, '55.00%'),
    @(95, 'Coder Camps', 11, 6, '54.55%'),
    @(96, 'Bloc.io', 21, 11, '52.38%'),
    @(97, 'Wyncode', 6, 3, '50.00%'),
    @(98, 'LearningFuze', 2, 1, '50.00%'),
    @(99, 'Fire Bootcamp', 2, 1, '50.00%'),
    @(100, 'HackerYou', 2, 1, '50.00%'),
    @(101, '4Geeks Academy', 8, 4, '50.00%'),
    @(102, 'Digital House', 2, 1, '50.00%'),
    @(103, 'Devschool', 4, 2, '50.00%'),
    @(104, 'Mobile Makers Academy', 2, 1, '50.00%'),
    @(105, 'Codesmith', 4, 2, '50.00%'),
    @(106, 'RefactorU', 4, 2, '50.00%'),
    @(107, 'Metis', 2, 1, '50.00%'),
    @(108, 'Bit Bootcamp', 4, 2, '50.00%'),
    @(109, 'PDX Code Guild', 2, 1, '50.00%'),
    @(110, 'Camp Code Away', 4, 2, '50.00%'),
    @(111, 'tradecraft', 3, 1, '33.33%'),
    @(112, 'Zip Code Wilmington', 3, 1, '33.33%'),
    @(113, 'Starter League', 3, 1, '33.33%'),
    @(114, 'Coding House', 3, 1, '33.33%'),
    @(115, 'We Can Code IT', 3, 1, '33.33%'),
    @(116, 'V School', 3, 1, '33.33%'),
    @(117, 'Codify Academy', 4, 1, '25.00%'),
    @(118, 'Code For Progress', 1, 0, '0.00%'),
    @(119, 'World Tech Makers', 1, 0, '0.00%'),
    @(120, 'DaVinci Coders', 2, 0, '0.00%'),
    @(121, 'Alphappl', 1, 0, '0.00%'),
    @(122, 'Academic Work Academy', 1, 0, '0.00%'),
    @(123, 'Astro Code School', 1, 0, '0.00%'),
    @(124, 'Atlanta Code', 1, 0, '0.00%'),
    @(125, 'Data Science Dojo', 1, 0, '0.00%'),
    @(126, 'Interface Web School', 1, 0, '0.00%'),
    @(127, 'Montana Code School', 1, 0, '0.00%'),
    @(128, 'TalentBuddy', 1, 0, '0.00%'),
    @(129, 'Academy X', 1, 0, '0.00%'),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
}

# Restore the default (unstyled) cell style on the % Recommend column so the
# output matches the original formatting (no explicit style index) while
# keeping the values stored as text.
$pctRange.Style = "Normal"

Write-Host "Updated $($data.Count) rows"